$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at row 41, shifting existing rows 41-66 down to 46-71
$ws.Rows("41:45").Insert()

$newData = New-Object 'object[,]' 5,20
$newData[0,0] = 6
$newData[0,1] = 'Mercado Mayorista Lo Valledor de Santiago'
$newData[0,2] = 'Metropolitana'
$newData[0,3] = 44539
$newData[0,4] = 13
$newData[0,5] = 'Fruta'
$newData[0,6] = 100103
$newData[0,7] = 'Frutos de hueso (carozo)'
$newData[0,8] = 100103003
$newData[0,9] = 'Damasco'
$newData[0,10] = 'Castle Brite'
$newData[0,11] = 'Especial'
$newData[0,12] = 100
$newData[0,13] = 18000
$newData[0,14] = 18000
$newData[0,15] = 18000
$newData[0,16] = '$/caja 15 kilos'
$newData[0,17] = 'Provincia de San Felipe de Aconcagua'
$newData[0,18] = 1200
$newData[0,19] = 15
$newData[1,0] = 6
$newData[1,1] = 'Mercado Mayorista Lo Valledor de Santiago'
$newData[1,2] = 'Metropolitana'
$newData[1,3] = 44539
$newData[1,4] = 13
$newData[1,5] = 'Fruta'
$newData[1,6] = 100103
$newData[1,7] = 'Frutos de hueso (carozo)'
$newData[1,8] = 100103003
$newData[1,9] = 'Damasco'
$newData[1,10] = 'Castle Brite'
$newData[1,11] = 'Especial'
$newData[1,12] = 65
$newData[1,13] = 19000
$newData[1,14] = 19000
$newData[1,15] = 19000
$newData[1,16] = '$/caja 18 kilos'
$newData[1,17] = 'Provincia de San Felipe de Aconcagua'
$newData[1,18] = 1056
$newData[1,19] = 18
$newData[2,0] = 6
$newData[2,1] = 'Mercado Mayorista Lo Valledor de Santiago'
$newData[2,2] = 'Metropolitana'
$newData[2,3] = 44539
$newData[2,4] = 13
$newData[2,5] = 'Fruta'
$newData[2,6] = 100103
$newData[2,7] = 'Frutos de hueso (carozo)'
$newData[2,8] = 100103003
$newData[2,9] = 'Damasco'
$newData[2,10] = 'Castle Brite'
$newData[2,11] = 'Primera'
$newData[2,12] = 185
$newData[2,13] = 16000
$newData[2,14] = 16000
$newData[2,15] = 16000
$newData[2,16] = '$/caja 15 kilos'
$newData[2,17] = 'Provincia de San Felipe de Aconcagua'
$newData[2,18] = 1067
$newData[2,19] = 15
$newData[3,0] = 6
$newData[3,1] = 'Mercado Mayorista Lo Valledor de Santiago'
$newData[3,2] = 'Metropolitana'
$newData[3,3] = 44539
$newData[3,4] = 13
$newData[3,5] = 'Fruta'
$newData[3,6] = 100103
$newData[3,7] = 'Frutos de hueso (carozo)'
$newData[3,8] = 100103003
$newData[3,9] = 'Damasco'
$newData[3,10] = 'Castle Brite'
$newData[3,11] = 'Primera'
$newData[3,12] = 220
$newData[3,13] = 15000
$newData[3,14] = 16000
$newData[3,15] = 15455
$newData[3,16] = '$/caja 18 kilos'
$newData[3,17] = 'Provincia de San Felipe de Aconcagua'
$newData[3,18] = 859
$newData[3,19] = 18
$newData[4,0] = 6
$newData[4,1] = 'Mercado Mayorista Lo Valledor de Santiago'
$newData[4,2] = 'Metropolitana'
$newData[4,3] = 44539
$newData[4,4] = 13
$newData[4,5] = 'Fruta'
$newData[4,6] = 100103
$newData[4,7] = 'Frutos de hueso (carozo)'
$newData[4,8] = 100103003
$newData[4,9] = 'Damasco'
$newData[4,10] = 'Castle Brite'
$newData[4,11] = 'Segunda'
$newData[4,12] = 175
$newData[4,13] = 12000
$newData[4,14] = 12000
$newData[4,15] = 12000
$newData[4,16] = '$/caja 15 kilos'
$newData[4,17] = 'Provincia de San Felipe de Aconcagua'
$newData[4,18] = 800
$newData[4,19] = 15

$ws.Range("A41:T45").Value = $newData
